# Egypt Division 1 - update bases das ligas (20-06-2024 20:11)
# The edit re-shuffles the per-match data (every column except the running
# index in column A) among several groups of rows, and refreshes a handful
# of odds/date values for two "upcoming match" rows at the bottom of the
# sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Groups of rows whose B:AD content rotates: row[i] ends up holding the
# data that row[i+1] (wrapping around) held before the edit.
$groups = @(
    @(23, 24),
    @(38, 39, 40),
    @(42, 43),
    @(107, 108),
    @(131, 132),
    @(137, 138),
    @(157, 158),
    @(172, 174),
    @(175, 177),
    @(185, 186),
    @(221, 222),
    @(230, 231),
    @(235, 236)
)

foreach ($group in $groups) {
    $count = $group.Count

    # Snapshot the B:AD values for every row in this group before writing
    # anything back (so later writes don't clobber data we still need).
    $snapshots = @()
    foreach ($row in $group) {
        $rng = $ws.Range("B$row`:AD$row")
        $snapshots += $rng.Value()
    }

    for ($i = 0; $i -lt $count; $i++) {
        $destRow = $group[$i]
        $srcIndex = ($i + 1) % $count
        $destRng = $ws.Range("B$destRow`:AD$destRow")
        $destRng.Value = $snapshots[$srcIndex]
    }
}

# Refresh closing odds (and kickoff time) for the two still-to-be-played
# fixtures at the end of the sheet.
$ws.Range("O255").Value = 1.7
$ws.Range("P255").Value = 3.4
$ws.Range("Q255").Value = 4.75
$ws.Range("S255").Value = 1.925
$ws.Range("T255").Value = 1.875
$ws.Range("V255").Value = 1.85
$ws.Range("W255").Value = 1.95

$ws.Range("D256").Value = 45464.58333333334
$ws.Range("O256").Value = 1.166
$ws.Range("P256").Value = 6
$ws.Range("S256").Value = 1.9
$ws.Range("T256").Value = 1.9
$ws.Range("V256").Value = 1.85
$ws.Range("W256").Value = 1.95
